$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to keep text representation (values like "53.70" or "1.002"
# would otherwise be auto-converted to numbers by Excel, dropping trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.897.44'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '1.701.79'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").Value = '314.83'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("D7").Value = '0.4007'
$ws.Range("E7").Value = '  +2.05%  '
$ws.Range("D8").Value = '0.4069'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = '1.004'
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").Value = '53.70'
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("D11").Value = '1.464'
$ws.Range("E11").Value = '  -2.89%  '
$ws.Range("D12").Value = '0.08814'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '25.90'
$ws.Range("E13").Value = '  +6.21%  '
$ws.Range("D14").Value = '7.482'
$ws.Range("E14").Value = '  -2.27%  '
$ws.Range("D15").Value = '8.040'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").Value = '0.00001346'
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").Value = '1.706.47'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '96.67'
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("D19").Value = '0.07192'
$ws.Range("E19").Value = '  +1.03%  '
$ws.Range("D20").Value = '20.95'
$ws.Range("E20").Value = '  +5.57%  '
$ws.Range("D21").Value = '7.247'
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").Value = '14.52'
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("D24").Value = '24.886.45'
$ws.Range("E24").Value = '  +0.97%  '
$ws.Range("D25").Value = '2.334'
$ws.Range("E25").Value = '  -1.19%  '
$ws.Range("D26").Value = '2.884'
$ws.Range("E26").Value = '  -5.39%  '
$ws.Range("D27").Value = '6.645'
$ws.Range("E27").Value = '  +26.59%  '
$ws.Range("D28").Value = '23.10'
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = '163.39'
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '143.66'
$ws.Range("E30").Value = '  +4.16%  '
$ws.Range("D31").Value = '8.203'
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("D32").Value = '1.938.62'
$ws.Range("E32").Value = '  +2.14%  '
$ws.Range("D33").Value = '2.271'
$ws.Range("E33").Value = '  +14.14%  '
$ws.Range("D34").Value = '0.08736'
$ws.Range("E34").Value = '  -1.30%  '
$ws.Range("D35").Value = '7.405'
$ws.Range("E35").Value = '  -1.51%  '
$ws.Range("D36").Value = '0.03167'
$ws.Range("E36").Value = '  +8.10%  '
$ws.Range("D37").Value = '1.032'
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").Value = '0.2871'
$ws.Range("E38").Value = '  +4.95%  '
$ws.Range("D39").Value = '0.8537'
$ws.Range("E39").Value = '  +8.17%  '
$ws.Range("D40").Value = '10.88'
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("D41").Value = '0.09440'
$ws.Range("E41").Value = '  +3.08%  '
$ws.Range("D42").Value = '14.07'
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("D43").Value = '1.474'
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("D44").Value = '17.62'
$ws.Range("E44").Value = '  +6.17%  '
$ws.Range("D45").Value = '2.696'
$ws.Range("E45").Value = '  +4.82%  '
$ws.Range("D46").Value = '0.7468'
$ws.Range("E46").Value = '  +3.34%  '
$ws.Range("D47").Value = '4.225'
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").Value = '1.404'
$ws.Range("E48").Value = '  +5.29%  '
$ws.Range("D49").Value = '1.004'
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").Value = '141.02'
$ws.Range("E50").Value = '  +1.12%  '
$ws.Range("D51").Value = '0.08358'
$ws.Range("E51").Value = '  +4.64%  '
